$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before F ("date" shifts from F to G)
$ws.Columns.Item(6).Insert()

# 2. Set the new F column header and bulk-fill population values for existing rows
$ws.Range("F1").Value = "population"
$ws.Range("F2:F361").Value = 218162

# 3. Insert 3 new rows after row 361 (rows 362-364)
$ws.Range("A362:A364").EntireRow.Insert()

# 4. Fix formatting of the newly inserted rows' column A to match the existing
#    style used throughout column A (style index 1, with border/bold/center/top)
$ws.Range("A361").Copy()
$ws.Range("A362:A364").PasteSpecial(-4122)

# 5. Update row 361 with corrected data values
$ws.Range("B361").Value = 10834
$ws.Range("C361").Value = 190
$ws.Range("D361").Value = 3
$ws.Range("E361").Value = 119

# 6. Fill in data for new row 362 (18/mar)
$ws.Range("A362").Value = "18/mar"
$ws.Range("B362").Value = 10834
$ws.Range("C362").Value = 190
$ws.Range("D362").Value = 0
$ws.Range("E362").Value = 0
$ws.Range("G362").Value = "18/mar"

# 7. Fill in data for new row 363 (19/mar)
$ws.Range("A363").Value = "19/mar"
$ws.Range("B363").Value = 11011
$ws.Range("C363").Value = 197
$ws.Range("D363").Value = 7
$ws.Range("E363").Value = 177
$ws.Range("G363").Value = "19/mar"

# 8. Fill in data for new row 364 (20/mar)
$ws.Range("A364").Value = "20/mar"
$ws.Range("B364").Value = 11011
$ws.Range("C364").Value = 197
$ws.Range("D364").Value = 0
$ws.Range("E364").Value = 0
$ws.Range("G364").Value = "20/mar"

# 9. Fill population column for the new rows as well
$ws.Range("F362:F364").Value = 218162
